{"js": "// Applies the content edits described by the diff to the \"Amazon RDS\",\n// \"Amazon Lambda\", \"Lambda (nachtelijke delete)\" and \"Opslag (S3 Bucket)\"\n// sections, plus two \"mysql\" -> \"MySQL\" capitalization fixes.\n\nconst body = context.document.body;\n\n// 1) \"...programmeertaal mysql om bestanden...\" -> \"...programmeertaal MySQL om bestanden...\"\nconst mysqlLower = body.search(\"mysql\", { matchCase: true });\nmysqlLower.load(\"items/text\");\nawait context.sync();\nfor (let i = 0; i < mysqlLower.items.length; i++) {\n  mysqlLower.items[i].insertText(\"MySQL\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Rewrite the \"We maken gebruik van Amazon Lambda...\" paragraph describing\n//    the automatic cleanup of old files, and delete the following paragraph\n//    (\"Hetzelfde moet gebeuren in onze Mysql AWS Relational Database...\")\n//    whose content was folded into the rewritten paragraph above.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet lambdaCleanupParagraph = null;\nlet duplicateRdsParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"We maken gebruik van Amazon Lambda om ervoor te zorgen dat\") >= 0) {\n    lambdaCleanupParagraph = paragraphs.items[i];\n  } else if (t.indexOf(\"Hetzelfde moet gebeuren in onze\") >= 0) {\n    duplicateRdsParagraph = paragraphs.items[i];\n  }\n}\n\nif (lambdaCleanupParagraph) {\n  lambdaCleanupParagraph.insertText(\n    \"We maken gebruik van Amazon Lambda om ervoor te zorgen dat alle bestanden die ouder zijn dan 24u dagelijks in onze MySQL AWS Relational Database om 3u in de nacht verwijderd worden. We willen namelijk, volgens de opdracht, dat bestanden niet langer dan 1 dag op ons systeem bewaard worden.\",\n    Word.InsertLocation.replace\n  );\n}\nif (duplicateRdsParagraph) {\n  duplicateRdsParagraph.delete();\n}\nawait context.sync();\n\n// 3) \"...verwijdert uit zowel de RDS als de S3 bucket.\" -> \"...verwijdert uit de RDS.\"\nconst rdsS3 = body.search(\"objecten dat ouder zijn dan 24u verwijdert uit zowel de RDS als de S3 bucket.\", { matchCase: true });\nrdsS3.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < rdsS3.items.length; i++) {\n  rdsS3.items[i].insertText(\"objecten dat ouder zijn dan 24u verwijdert uit de RDS.\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 4) Append a sentence about the S3 bucket's built-in lifecycle rule to the\n//    end of the \"Opslag (S3 Bucket)\" paragraph.\nconst s3Paragraphs = body.paragraphs;\ns3Paragraphs.load(\"items/text\");\nawait context.sync();\n\nlet s3BucketParagraph = null;\nfor (let i = 0; i < s3Paragraphs.items.length; i++) {\n  const t = s3Paragraphs.items[i].text;\n  if (t.indexOf(\"De S3 Bucket is een eenvoudige opslag service\") >= 0) {\n    s3BucketParagraph = s3Paragraphs.items[i];\n  }\n}\nif (s3BucketParagraph) {\n  s3BucketParagraph.insertText(\n    \" De bestanden worden verwijderd als ze ouder zijn dan 24u met behulp van de ingebouwde service life cycle.\",\n    Word.InsertLocation.end\n  );\n}\nawait context.sync();\n\n// 5) \"We hebben ook gebruik gemaakt van Mysql...\" -> \"...MySQL...\"\nconst mysqlUpper = body.search(\"Mysql\", { matchCase: true });\nmysqlUpper.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < mysqlUpper.items.length; i++) {\n  mysqlUpper.items[i].insertText(\"MySQL\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Applies the content edits described by the diff to the \"Amazon RDS\",\n# \"Amazon Lambda\", \"Lambda (nachtelijke delete)\" and \"Opslag (S3 Bucket)\"\n# sections, plus two \"mysql\" -> \"MySQL\" capitalization fixes.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.MatchCase = $true\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\n# 1) \"...programmeertaal mysql om bestanden...\" -> \"...programmeertaal MySQL om bestanden...\"\nReplace-Text \"mysql\" \"MySQL\"\n\n# 2) Rewrite the \"We maken gebruik van Amazon Lambda...\" paragraph describing\n#    the automatic cleanup of old files.\nReplace-Text \"alle objecten (bestanden) binnen onze Amazon S3 bucket automatisch na 24 uur verwijderd worden. We willen namelijk, volgens de opdracht, dat bestanden niet langer dan 1 op ons systeem bewaard worden.\" \"alle bestanden die ouder zijn dan 24u dagelijks in onze MySQL AWS Relational Database om 3u in de nacht verwijderd worden. We willen namelijk, volgens de opdracht, dat bestanden niet langer dan 1 dag op ons systeem bewaard worden.\"\n\n# ...and delete the following paragraph (\"Hetzelfde moet gebeuren in onze Mysql\n# AWS Relational Database...\") whose content was folded into the rewritten\n# paragraph above.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Hetzelfde moet gebeuren in onze*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 3) \"...verwijdert uit zowel de RDS als de S3 bucket.\" -> \"...verwijdert uit de RDS.\"\nReplace-Text \"objecten dat ouder zijn dan 24u verwijdert uit zowel de RDS als de S3 bucket.\" \"objecten dat ouder zijn dan 24u verwijdert uit de RDS.\"\n\n# 4) Append a sentence about the S3 bucket's built-in lifecycle rule to the\n#    end of the \"Opslag (S3 Bucket)\" paragraph.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*De S3 Bucket is een eenvoudige opslag service*\") {\n        $p.Range.InsertAfter(\" De bestanden worden verwijderd als ze ouder zijn dan 24u met behulp van de ingebouwde service life cycle.\")\n        break\n    }\n}\n\n# 5) \"We hebben ook gebruik gemaakt van Mysql...\" -> \"...MySQL...\"\nReplace-Text \"Mysql\" \"MySQL\"\n"}
